$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input_Value")

# Row 2 held several "AutomationTesting_3" test-data values; bump them to
# "AutomationTesting_4" (the "Accrual Clearing" value in L2 is untouched).
$ws.Range("I2").Value = "AutomationTesting_4"
$ws.Range("J2").Value = "AutomationTesting_4"
$ws.Range("M2").Value = "AutomationTesting_4"
$ws.Range("N2").Value = "AutomationTesting_4"

# Leave the cursor parked on H12, matching the saved selection.
$ws.Activate()
$ws.Range("H12").Select()
